# Wafer map DB load success.
# Fills in the previously-placeholder Wafer/Sub IDs (column D, rows 2 and
# 4-25) with the values returned by the DB load, applies the thin-border
# "loaded" cell format to those cells, repositions the window, and updates
# the active selection - matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Sub ID values coming back from the DB load (row 3 was not part of
# this batch and keeps its original placeholder value).
$loaded = @{
    2  = "86809619CNG7"
    4  = "65420204CNF0"
    5  = "63518033CNC5"
    6  = "63617427CNF6"
    7  = "65220505CNC7"
    8  = "18009303CNC5"
    9  = "75119936CNC7"
    10 = "79817714CNA4"
    11 = "71319718CNG7"
    12 = "61016332CNH2"
    13 = "73716819CNC3"
    14 = "66316625CNF0"
    15 = "84207528CNC7"
    16 = "62317612CNB7"
    17 = "70819903CND2"
    18 = "75320025CNF3"
    19 = "62117517CNF4"
    20 = "75719206CNG6"
    21 = "72119624CNG4"
    22 = "73818504CNE0"
    23 = "82408825CNF7"
    24 = "63617406CNH2"
    25 = "63617408CNB0"
}

foreach ($row in 2..25) {
    if ($loaded.ContainsKey($row)) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.Value = $loaded[$row]
        $cell.Borders.LineStyle = 1
    }
}

# Move the active selection (matches the recorded cursor position after the load).
[void]$ws.Range("M10").Select()

# Reposition / resize the document window (best-effort - matches recorded window geometry).
$win = $excel.ActiveWindow
$win.Left = 4642
$win.Top = 3555
$win.Width = 21601
$win.Height = 12683
